# Apply "Added for AM work" changes to the TRIP780F.xpc sheet.
# 1. Row 14 (HexGrid-90degTilt5degRes) is renamed and its data values are refreshed.
# 2. Two new rows (24, 25) are appended for the new "RotRing Axis-Y" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row 14 - label and values
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "HexGrid-90degTilt22p5degRes"

$ws.Range("C14").Value = 0.9942783707768669
$ws.Range("D14").Value = 0.9876431694675819
$ws.Range("E14").Value = 0.9801133740497162
$ws.Range("F14").Value = 0.9942783707768669
$ws.Range("G14").Value = 0.9881620024722707
$ws.Range("H14").Value = 0.9631200627381333
$ws.Range("I14").Value = 0.9859554117093244
$ws.Range("J14").Value = 0.9876431694675819
$ws.Range("K14").Value = 0.9838782717586491
$ws.Range("L14").Value = 0.9838782717586491
$ws.Range("M14").Value = 0.9853061819965229
$ws.Range("N14").Value = 0.9873449714313883
$ws.Range("O14").Value = 0.9873449714313883
$ws.Range("P14").Value = 0.989078321267758
$ws.Range("Q14").Value = 0.989078321267758
$ws.Range("R14").Value = 0.9832120652023155

# ---------------------------------------------------------------------------
# 2. Append new row 24 - "RotRing Axis-Y Res-5.0 Theta-2.5 "
# ---------------------------------------------------------------------------
# Copy formatting from the preceding row's A cell (bold / centered / bordered)
# so the new index cells match the look of the rest of column A.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "RotRing Axis-Y Res-5.0 Theta-2.5 "
$ws.Range("C24").Value = 1.222691003359202
$ws.Range("D24").Value = 0.895576178072181
$ws.Range("E24").Value = 0.9709247760303786
$ws.Range("F24").Value = 1.222691003359202
$ws.Range("G24").Value = 0.9375911349841114
$ws.Range("H24").Value = 0.9196952270909058
$ws.Range("I24").Value = 0.9983055836563478
$ws.Range("J24").Value = 0.895576178072181
$ws.Range("K24").Value = 0.9332504770512798
$ws.Range("L24").Value = 0.9332504770512798
$ws.Range("M24").Value = 0.9346973630288904
$ws.Range("N24").Value = 1.029730652487254
$ws.Range("O24").Value = 1.029730652487254
$ws.Range("P24").Value = 1.077970740205241
$ws.Range("Q24").Value = 1.077970740205241
$ws.Range("R24").Value = 0.9907973171988544

# ---------------------------------------------------------------------------
# 3. Append new row 25 - "RotRing Axis-Y Res-5.0 Theta-5.0 "
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "RotRing Axis-Y Res-5.0 Theta-5.0 "
$ws.Range("C25").Value = 1.145805252063063
$ws.Range("D25").Value = 0.9003833225384833
$ws.Range("E25").Value = 0.9719831634697815
$ws.Range("F25").Value = 1.145805252063063
$ws.Range("G25").Value = 0.9438708744676136
$ws.Range("H25").Value = 0.9255320281711241
$ws.Range("I25").Value = 0.9990716574632709
$ws.Range("J25").Value = 0.9003833225384833
$ws.Range("K25").Value = 0.9361832430041324
$ws.Range("L25").Value = 0.9361832430041324
$ws.Range("M25").Value = 0.9387457868252929
$ws.Range("N25").Value = 1.006057246023776
$ws.Range("O25").Value = 1.006057246023776
$ws.Range("P25").Value = 1.040994247533598
$ws.Range("Q25").Value = 1.040994247533598
$ws.Range("R25").Value = 0.9811077163622227
